$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period header labels in row 2-4 (C column) to shift from */12 to */03
$ws.Range("C2").Value = "2015/03  (IFRS연결)"
$ws.Range("C3").Value = "2016/03  (IFRS연결)"
$ws.Range("C4").Value = "2017/03  (IFRS연결)"

# Row 2: replace financial figures with restated (smaller-scale) values
$ws.Range("D2").Value = 628
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = -77
$ws.Range("H2").Value = -74
$ws.Range("I2").Value = -75
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1014
$ws.Range("L2").Value = 299
$ws.Range("M2").Value = 715
$ws.Range("N2").Value = 704
$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 77
$ws.Range("Q2").Value = 14
$ws.Range("R2").Value = 29
$ws.Range("S2").Value = -5
$ws.Range("T2").Value = 4
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 50
$ws.Range("W2").Value = 2.18
$ws.Range("X2").Value = -11.72
$ws.Range("Y2").Value = -10.53
$ws.Range("Z2").Value = -7.45
$ws.Range("AA2").Value = 41.84
$ws.Range("AB2").Value = 764.7
$ws.Range("AC2").Value = -1035
$ws.Range("AD2").Value = -9.710000000000001
$ws.Range("AE2").Value = 10659
$ws.Range("AF2").Value = 0.9399999999999999
$ws.Range("AG2").Value = 180
$ws.Range("AH2").Value = 1.79
$ws.Range("AI2").Value = -16.1
$ws.Range("AJ2").Value = 5984918

# Row 3: replace financial figures with restated (smaller-scale) values
$ws.Range("D3").Value = 697
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 16
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = -4
$ws.Range("I3").Value = -3
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 1043
$ws.Range("L3").Value = 306
$ws.Range("M3").Value = 737
$ws.Range("N3").Value = 729
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 78
$ws.Range("Q3").Value = -117
$ws.Range("R3").Value = -35
$ws.Range("S3").Value = 139
$ws.Range("T3").Value = 21
$ws.Range("U3").Value = -138
$ws.Range("V3").Value = 174
$ws.Range("W3").Value = 2.31
$ws.Range("X3").Value = -0.58
$ws.Range("Y3").Value = -0.35
$ws.Range("Z3").Value = -0.39
$ws.Range("AA3").Value = 41.5
$ws.Range("AB3").Value = 753.88
$ws.Range("AC3").Value = -34
$ws.Range("AD3").Value = -390.44
$ws.Range("AE3").Value = 10399
$ws.Range("AF3").Value = 1.28
$ws.Range("AG3").Value = 180
$ws.Range("AH3").Value = 1.35
$ws.Range("AI3").Value = -504.24
$ws.Range("AJ3").Value = 6161406

# Row 4: replace financial figures with restated (smaller-scale) values
$ws.Range("D4").Value = 716
$ws.Range("E4").Value = 37
$ws.Range("F4").Value = 37
$ws.Range("G4").Value = 65
$ws.Range("H4").Value = 61
$ws.Range("I4").Value = 59
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1113
$ws.Range("L4").Value = 330
$ws.Range("M4").Value = 783
$ws.Range("N4").Value = 772
$ws.Range("O4").Value = 11
$ws.Range("P4").Value = 78
$ws.Range("Q4").Value = 98
$ws.Range("R4").Value = 15
$ws.Range("S4").Value = -22
$ws.Range("T4").Value = 22
$ws.Range("U4").Value = 76
$ws.Range("V4").Value = 169
$ws.Range("W4").Value = 5.22
$ws.Range("X4").Value = 8.59
$ws.Range("Y4").Value = 7.85
$ws.Range("Z4").Value = 5.7
$ws.Range("AA4").Value = 42.16
$ws.Range("AB4").Value = 809.16
$ws.Range("AC4").Value = 789
$ws.Range("AD4").Value = 11.38
$ws.Range("AE4").Value = 11013
$ws.Range("AF4").Value = 0.82
$ws.Range("AG4").Value = 180
$ws.Range("AH4").Value = 2
$ws.Range("AI4").Value = 21.65
$ws.Range("AJ4").Value = 6161406

# Row 5: replace financial figures with restated (smaller-scale) values
$ws.Range("D5").Value = 629
$ws.Range("E5").Value = 36
$ws.Range("F5").Value = 36
$ws.Range("G5").Value = 48
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 45
$ws.Range("J5").Value = -3
$ws.Range("K5").Value = 1215
$ws.Range("L5").Value = 387
$ws.Range("M5").Value = 828
$ws.Range("N5").Value = 820
$ws.Range("O5").Value = 9
$ws.Range("P5").Value = 80
$ws.Range("Q5").Value = 13
$ws.Range("R5").Value = -27
$ws.Range("S5").Value = 23
$ws.Range("T5").Value = 45
$ws.Range("U5").Value = -32
$ws.Range("V5").Value = 194
$ws.Range("W5").Value = 5.68
$ws.Range("X5").Value = 6.82
$ws.Range("Y5").Value = 5.71
$ws.Range("Z5").Value = 3.68
$ws.Range("AA5").Value = 46.67
$ws.Range("AB5").Value = 839
$ws.Range("AC5").Value = 607
$ws.Range("AD5").Value = 19.61
$ws.Range("AE5").Value = 11362
$ws.Range("AF5").Value = 1.05
$ws.Range("AG5").Value = 190
$ws.Range("AH5").Value = 1.6
$ws.Range("AI5").Value = 30.47
$ws.Range("AJ5").Value = 6368587

# Row 6: replace financial figures with restated (smaller-scale) values
$ws.Range("D6").Value = 831
$ws.Range("E6").Value = 53
$ws.Range("F6").Value = 53
$ws.Range("G6").Value = 46
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 40
$ws.Range("K6").Value = 1352
$ws.Range("L6").Value = 484
$ws.Range("M6").Value = 869
$ws.Range("N6").Value = 858
$ws.Range("P6").Value = 80
$ws.Range("Q6").Value = 61
$ws.Range("R6").Value = -99
$ws.Range("S6").Value = 109
$ws.Range("T6").Value = 98
$ws.Range("U6").Value = -37
$ws.Range("V6").Value = 272
$ws.Range("W6").Value = 6.41
$ws.Range("X6").Value = 5.19
$ws.Range("Y6").Value = 4.79
$ws.Range("Z6").Value = 3.36
$ws.Range("AA6").Value = 55.66
$ws.Range("AB6").Value = 931.9299999999999
$ws.Range("AC6").Value = 524
$ws.Range("AD6").Value = 20.02
$ws.Range("AE6").Value = 11991
$ws.Range("AF6").Value = 0.88
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 1.9
$ws.Range("AI6").Value = 35.91
$ws.Range("AJ6").Value = 6368587

# Rows 7-9: clear all data columns (D:AJ), keep only A/B/C identifying columns
$ws.Range("D7:AJ9").ClearContents()
